# "Packages können erstellt un aquired werden" -
# mark the "Create and acquire packages" feature row as checked/done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D next to "Create and acquire packages" (row 39) and
# "Show and configure decks" (row 40) gets a "check" comment/mark.
$ws.Range("D39").Value = "check"
$ws.Range("D40").Value = "check"

# "Contains link to GIT" (row 51) was missing its awarded points - fill it in.
$ws.Range("C51").Value = 0.5

# Reflect where the author was scrolled to / had selected when saving.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D23").Select()
